# Generate Report for Archive
#
# 1. Status text "Ready for handoff" -> "In Translation" everywhere it is
#    used (Overview!E2:F3, zh-cn!C2:C3, de-de!C2:C3 all shared the same
#    string).
# 2. Narrow the "Status" columns (Overview E & F, zh-cn C, de-de C).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Update the status text wherever it appears ---
$overview.Range("E2:F2").Value = "In Translation"
$overview.Range("E3:F3").Value = "In Translation"

$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# --- Narrow the Status columns ---
$overview.Columns("E:F").ColumnWidth = 12.5
$zhcn.Columns("C:C").ColumnWidth = 12.5
$dede.Columns("C:C").ColumnWidth = 12.5
